$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Number" -> "Label" for the new Label column header
$ws.Range("C2").Value = "Label"

# Opamp part number correction
$ws.Range("A6").Value = "LM324D opamps"

# New "Label" column values (which component reference designators use each part)
$ws.Range("C3").Value = "C1, C2, C3, C4"
$ws.Range("C4").Value = "R6, R11, R12"
$ws.Range("C5").Value = "R1, R2, R3, R4, R7, R8, R9, R10"
$ws.Range("C6").Value = "U9, U10"
$ws.Range("C7").Value = "U3, U4, U7, U8"
$ws.Range("C8").Value = "U1, U2, U5, U6"
$ws.Range("C9").Value = "DAC1"
$ws.Range("C10").Value = "R5"
$ws.Range("C11").Value = "R13, R14"

# Match the selected cell left by the editor
$ws.Range("E8").Select() | Out-Null
